$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Name" / "Amount Spent" mini-table in columns F:G, mirroring the
# bold+underline+centered header style already used by A1:D1.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "Name"
$ws.Range("G1").Value = "Amount Spent"
$ws.Range("F2").Value = "Krushik"
$ws.Range("F3").Value = "Likith"
$ws.Range("F4").Value = "Nandan"
$ws.Range("G4").Value = 4934
$ws.Range("F5").Value = "Nikhil"

# Widen the two new columns to fit their content.
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 14

# Leave the selection where the author last clicked.
$ws.Range("H4").Select()
